# Update cryptos list: refresh Price (column D) and Volume(1h) (column E)
# values for the coin rows that changed, matching the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.571.15"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "3.479.17"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("D5").Value = "'578.70"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").Value = "'160.30"
$ws.Range("E6").Value = "  +2.46%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.478.38"
$ws.Range("E8").Value = "  +0.35%  "
$ws.Range("D9").Value = "'0.583"
$ws.Range("E9").Value = "  +4.45%  "
$ws.Range("D10").Value = "'7.26"
$ws.Range("E10").Value = "  -4.13%  "
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("D12").Value = "'0.441"
$ws.Range("E12").Value = "  -1.59%  "
$ws.Range("D13").Value = "4.073.75"
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("E14").Value = "  -1.53%  "
$ws.Range("E15").Value = "  -2.78%  "
$ws.Range("D16").Value = "'28.68"
$ws.Range("E16").Value = "  +2.58%  "
$ws.Range("D17").Value = "65.554.47"
$ws.Range("E17").Value = "  +1.21%  "
$ws.Range("D18").Value = "3.454.32"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").Value = "'6.41"
$ws.Range("E19").Value = "  -0.66%  "
$ws.Range("D20").Value = "'14.26"
$ws.Range("E20").Value = "  -1.11%  "
$ws.Range("D21").Value = "'390.07"
$ws.Range("E21").Value = "  -2.03%  "
$ws.Range("D22").Value = "'8.24"
$ws.Range("E22").Value = "  -3.84%  "
$ws.Range("D23").Value = "'0.549"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "'73.52"
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").Value = "'9.57"
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("D28").Value = "'0.178"
$ws.Range("E28").Value = "  -1.28%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").Value = "'6.40"
$ws.Range("E30").Value = "  +5.62%  "
$ws.Range("E31").Value = "  +2.87%  "
$ws.Range("E32").Value = "  +0.49%  "
$ws.Range("D33").Value = "'23.68"
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("D34").Value = "'6.46"
$ws.Range("E34").Value = "  -4.92%  "
$ws.Range("D36").Value = "'7.08"
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("E37").Value = "  +2.71%  "
$ws.Range("D38").Value = "'162.70"
$ws.Range("E38").Value = "  +1.05%  "
$ws.Range("E39").Value = "  +3.58%  "
$ws.Range("D40").Value = "3.059.24"
$ws.Range("E40").Value = "  +5.12%  "
$ws.Range("D41").Value = "'0.0769"
$ws.Range("E41").Value = "  -2.06%  "
$ws.Range("D42").Value = "'27.03"
$ws.Range("E42").Value = "  -2.72%  "
$ws.Range("D43").Value = "'0.0319"
$ws.Range("E43").Value = "  -1.76%  "
$ws.Range("E44").Value = "  +1.37%  "
$ws.Range("D45").Value = "'42.72"
$ws.Range("E45").Value = "  +1.40%  "
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").Value = "'25.69"
$ws.Range("E47").Value = "  +7.98%  "
$ws.Range("E48").Value = "  +0.94%  "
$ws.Range("D49").Value = "'2.21"
$ws.Range("E49").Value = "  +0.72%  "
$ws.Range("E50").Value = "  +1.74%  "
$ws.Range("D51").Value = "'310.19"
$ws.Range("E51").Value = "  +3.51%  "
